# Update "想去人数" (F column) values on the "展览" (sheet1) and
# "全部类型" (sheet4) worksheets to match the regenerated site data.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# row => new F value, for the "展览" sheet
$exhibitionUpdates = @{
    2  = 3057
    3  = 475
    4  = 55
    5  = 40
    6  = 262
    7  = 1043
    8  = 14733
    9  = 174
    11 = 5875
    12 = 600
    16 = 1242
    18 = 92
    19 = 191
    20 = 806
    23 = 10676
    24 = 1206
    25 = 71
    26 = 105
    27 = 3746
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# row => new F value, for the "全部类型" sheet
$allTypesUpdates = @{
    3  = 3057
    4  = 475
    5  = 55
    6  = 40
    7  = 262
    8  = 1043
    9  = 14733
    10 = 174
    12 = 5875
    13 = 600
    17 = 1242
    19 = 92
    20 = 191
    21 = 806
    25 = 10676
    26 = 1206
    27 = 71
    28 = 105
    29 = 3746
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
